# "New crime data collected" -- weekly refresh of the CompStat_1 report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -----------------------------------------------
# A8 holds "Volume 31   Number  12" (rich text); bump the issue number 12 -> 13.
$volCell = $ws.Range("A8")
$volText = $volCell.Text
$volLen = $volText.Length
$numChars = $volCell.Characters($volLen - 1, 2)
$numChars.Text = "13"

# C9 holds "Report Covering the Week  3/18/2024  Through  3/24/2024"; roll
# the reporting week forward by one week (new week-ending data was collected).
$weekCell = $ws.Range("C9")
$weekText = $weekCell.Text
$startIdx = $weekText.IndexOf("3/18/2024") + 1
$startChars = $weekCell.Characters($startIdx, 9)
$startChars.Text = "3/25/2024"

$weekText2 = $weekCell.Text
$endIdx = $weekText2.IndexOf("3/24/2024") + 1
$endChars = $weekCell.Characters($endIdx, 9)
$endChars.Text = "3/31/2024"

# --- Weekly crime statistics table (rows 14-33) -------------------------
# Murder
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = -100
$ws.Range("J14").Value = 4

# Rape
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -60
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = -14.285714285714
$ws.Range("L15").Value = 200
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -71.428571428571

# Robbery
$ws.Range("C16").Value = 9
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 16.666666666666
$ws.Range("I16").Value = 33
$ws.Range("J16").Value = 52
$ws.Range("K16").Value = -36.538461538461
$ws.Range("L16").Value = 57.142857142857
$ws.Range("M16").Value = -50.746268656716
$ws.Range("N16").Value = -87.209302325581

# Fel. Assault
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = 26.666666666666
$ws.Range("I17").Value = 122
$ws.Range("J17").Value = 89
$ws.Range("K17").Value = 37.078651685393
$ws.Range("L17").Value = 40.229885057471
$ws.Range("M17").Value = 67.123287671232
$ws.Range("N17").Value = -27.380952380952

# Burglary
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = -68.75
$ws.Range("I18").Value = 24
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = -42.857142857142
$ws.Range("L18").Value = 33.333333333333
$ws.Range("M18").Value = -57.142857142857
$ws.Range("N18").Value = -94.326241134751

# Gr. Larceny
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 3.333333333333
$ws.Range("I19").Value = 90
$ws.Range("J19").Value = 105
$ws.Range("K19").Value = -14.285714285714
$ws.Range("L19").Value = -3.225806451612
$ws.Range("M19").Value = 16.883116883116
$ws.Range("N19").Value = -25.619834710743

# G.L.A.
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 71.428571428571
$ws.Range("I20").Value = 27
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = -10
$ws.Range("M20").Value = -38.636363636363
$ws.Range("N20").Value = -92.458100558659

# TOTAL
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -4.347826086956
$ws.Range("F21").Value = 102
$ws.Range("G21").Value = 101
$ws.Range("H21").Value = 0.990099009900
$ws.Range("I21").Value = 302
$ws.Range("J21").Value = 329
$ws.Range("K21").Value = -8.206686930091
$ws.Range("L21").Value = 25.833333333333
$ws.Range("M21").Value = -6.790123456790
$ws.Range("N21").Value = -77.695716395864

# Transit (row 22) - no numeric changes, stays all "***.*"

# Housing
$ws.Range("C23").Value = 3
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 12
$ws.Range("H23").Value = 140
$ws.Range("I23").Value = 21
$ws.Range("J23").Value = 25
$ws.Range("K23").Value = -16
$ws.Range("L23").Value = 50
$ws.Range("M23").Value = 162.5

# Petit Larceny
$ws.Range("C24").Value = 34
$ws.Range("E24").Value = 61.904761904761
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = 32.967032967033
$ws.Range("I24").Value = 351
$ws.Range("J24").Value = 294
$ws.Range("K24").Value = 19.387755102040
$ws.Range("L24").Value = 45.643153526971
$ws.Range("M24").Value = 12.5

# Retail Theft
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 114.285714285714
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 100
$ws.Range("I25").Value = 169
$ws.Range("J25").Value = 104
$ws.Range("K25").Value = 62.5
$ws.Range("L25").Value = 119.480519480519

# Misd. Assault
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 18
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 43
$ws.Range("G26").Value = 50
$ws.Range("H26").Value = -14
$ws.Range("I26").Value = 150
$ws.Range("J26").Value = 169
$ws.Range("K26").Value = -11.242603550295
$ws.Range("L26").Value = 6.382978723404
$ws.Range("M26").Value = -46.428571428571

# UCR Rape*
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -57.142857142857
$ws.Range("I27").Value = 12
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = 100

# Other Sex Crimes
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 23
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 43.75
$ws.Range("L28").Value = 76.923076923076

# Shooting Vic.
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = "0"
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = -100
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = 100
$ws.Range("M29").Value = -55.555555555555
$ws.Range("N29").Value = -77.777777777777

# Shooting Inc.
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = "0"
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = -100
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = 100
$ws.Range("M30").Value = -50
$ws.Range("N30").Value = -75

# Hate Crimes (row 31) - no numeric changes, stays all "***.*"

# Traffic Fatalities
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = -100
$ws.Range("F33").Value = "0"
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = -100
$ws.Range("J33").Value = 1
$ws.Range("K33").Value = 0
